# Update the dSF column (column F) values for a set of rows, per the
# "repull data, push all data, mean calculation" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    6  = -2
    16 = 0
    18 = 2
    22 = -4
    23 = -5
    26 = -2
    31 = 2
    38 = 0
    42 = -3
    44 = -7
    46 = 0
    47 = 1
    53 = -2
    55 = 0
    69 = -3
    70 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
